$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update computed result values for rows 2-25 (case with 380 kV)
$ws.Range("B2").Value = 1.436810351341705
$ws.Range("C2").Value = 0.1490456417542134
$ws.Range("D2").Value = 0.1011546772977816
$ws.Range("F2").Value = 1.976500301083391
$ws.Range("G2").Value = 1.353208153423566
$ws.Range("H2").Value = 1.258213683796825
$ws.Range("I2").Value = 1.170605710526154
$ws.Range("J2").Value = 0.1642800361228494
$ws.Range("L2").Value = 0.3747866254339414
$ws.Range("M2").Value = 0.365725907851143
$ws.Range("N2").Value = 1.748690795975323

$ws.Range("B3").Value = 1.352828245702085
$ws.Range("C3").Value = 0.1307920858440923
$ws.Range("D3").Value = 0.1010972954726412
$ws.Range("F3").Value = 1.977121544079765
$ws.Range("G3").Value = 1.349168409930797
$ws.Range("H3").Value = 1.261944361666991
$ws.Range("I3").Value = 1.177221593430588
$ws.Range("J3").Value = 0.1648917161391541
$ws.Range("L3").Value = 0.3722787096958911
$ws.Range("M3").Value = 0.3525175906403462
$ws.Range("N3").Value = 1.76799718991782

$ws.Range("B4").Value = 1.301780066469689
$ws.Range("C4").Value = 0.1195368416412634
$ws.Range("D4").Value = 0.1010730721646311
$ws.Range("F4").Value = 1.97855682623495
$ws.Range("G4").Value = 1.347545703292269
$ws.Range("H4").Value = 1.264837325929662
$ws.Range("I4").Value = 1.181888293977849
$ws.Range("J4").Value = 0.1652888628122842
$ws.Range("L4").Value = 0.3708780484156975
$ws.Range("M4").Value = 0.3445547826051936
$ws.Range("N4").Value = 1.78046312095403

$ws.Range("B5").Value = 1.281108580751265
$ws.Range("C5").Value = 0.1149383217560285
$ws.Range("D5").Value = 0.1010659810274372
$ws.Range("F5").Value = 1.979406643612265
$ws.Range("G5").Value = 1.347099846580278
$ws.Range("H5").Value = 1.266167694111502
$ws.Range("I5").Value = 1.183941943791766
$ws.Range("J5").Value = 0.1654561409253912
$ws.Range("L5").Value = 0.3703423709467515
$ws.Range("M5").Value = 0.3413470699169281
$ws.Range("N5").Value = 1.785696903464885

$ws.Range("B6").Value = 1.277684042112952
$ws.Range("C6").Value = 0.1141740220613485
$ws.Range("D6").Value = 0.1010649718037957
$ws.Range("F6").Value = 1.979563754252489
$ws.Range("G6").Value = 1.347038814783104
$ws.Range("H6").Value = 1.266397748849585
$ws.Range("I6").Value = 1.184292124740246
$ws.Range("J6").Value = 0.1654842461298331
$ws.Range("L6").Value = 0.3702555451661738
$ws.Range("M6").Value = 0.3408166839553388
$ws.Range("N6").Value = 1.786575258876741

$ws.Range("B7").Value = 1.301500751573911
$ws.Range("C7").Value = 0.1194748726095156
$ws.Range("D7").Value = 0.1010729652580036
$ws.Range("F7").Value = 1.978567214578192
$ws.Range("G7").Value = 1.347538818487365
$ws.Range("H7").Value = 1.264854654473766
$ws.Range("I7").Value = 1.181915375162433
$ws.Range("J7").Value = 0.1652910967494856
$ws.Range("L7").Value = 0.370870681818765
$ws.Range("M7").Value = 0.3445113714068171
$ws.Range("N7").Value = 1.780533082801428

$ws.Range("B8").Value = 1.407746560700048
$ws.Range("C8").Value = 0.1427616809324093
$ws.Range("D8").Value = 0.1011326139220188
$ws.Range("F8").Value = 1.976495735461569
$ws.Range("G8").Value = 1.351637064010021
$ws.Range("H8").Value = 1.25937501213059
$ws.Range("I8").Value = 1.172761348667422
$ws.Range("J8").Value = 0.1644864738256224
$ws.Range("L8").Value = 0.3738930604291397
$ws.Range("M8").Value = 0.3611412497157644
$ws.Range("N8").Value = 1.755220672238369

$ws.Range("B9").Value = 1.620164978638456
$ws.Range("C9").Value = 0.1880513564979367
$ws.Range("D9").Value = 0.1013364658528104
$ws.Range("F9").Value = 1.980800926212325
$ws.Range("G9").Value = 1.366493908784378
$ws.Range("H9").Value = 1.25340940911444
$ws.Range("I9").Value = 1.159611239664514
$ws.Range("J9").Value = 0.1630791863045165
$ws.Range("L9").Value = 0.3809208415273559
$ws.Range("M9").Value = 0.3949134766709008
$ws.Range("N9").Value = 1.710436654721656

$ws.Range("B10").Value = 1.77868362146404
$ws.Range("C10").Value = 0.2211013667765371
$ws.Range("D10").Value = 0.1015386224269044
$ws.Range("F10").Value = 1.989076462937504
$ws.Range("G10").Value = 1.381590862359275
$ws.Range("H10").Value = 1.251943297711279
$ws.Range("I10").Value = 1.152883546072587
$ws.Range("J10").Value = 0.162148407152281
$ws.Range("L10").Value = 0.3867513461017751
$ws.Range("M10").Value = 0.4204279896953338
$ws.Range("N10").Value = 1.680491412625923

$ws.Range("B11").Value = 1.851326366422768
$ws.Range("C11").Value = 0.2360891689713185
$ws.Range("D11").Value = 0.1016418476433216
$ws.Range("F11").Value = 1.993954129656117
$ws.Range("G11").Value = 1.38937224894616
$ws.Range("H11").Value = 1.251910437271761
$ws.Range("I11").Value = 1.150461557541128
$ws.Range("J11").Value = 0.1617471938973596
$ws.Range("L11").Value = 0.3895479017637911
$ws.Range("M11").Value = 0.4321866228945979
$ws.Range("N11").Value = 1.667510530937598

$ws.Range("B12").Value = 1.878909997458436
$ws.Range("C12").Value = 0.2417579526631926
$ws.Range("D12").Value = 0.1016825456499397
$ws.Range("F12").Value = 1.995961399066829
$ws.Range("G12").Value = 1.392450618997231
$ws.Range("H12").Value = 1.251989213852482
$ws.Range("I12").Value = 1.14963634590746
$ws.Range("J12").Value = 0.1615984443546452
$ws.Range("L12").Value = 0.3906275421997663
$ws.Range("M12").Value = 0.4366610124562342
$ws.Range("N12").Value = 1.662687265345848

$ws.Range("B13").Value = 1.872966032601255
$ws.Range("C13").Value = 0.2405373804758142
$ws.Range("D13").Value = 0.101673709221938
$ws.Range("F13").Value = 1.995521970320794
$ws.Range("G13").Value = 1.391781772900515
$ws.Range("H13").Value = 1.251968190211102
$ws.Range("I13").Value = 1.149809978710358
$ws.Range("J13").Value = 0.1616303389273357
$ws.Range("L13").Value = 0.390394105600933
$ws.Range("M13").Value = 0.4356964128803611
$ws.Range("N13").Value = 1.663721934843302

$ws.Range("B14").Value = 1.853594183357757
$ws.Range("C14").Value = 0.2365556791610004
$ws.Range("D14").Value = 0.1016451637099891
$ws.Range("F14").Value = 1.994116057388936
$ws.Range("G14").Value = 1.389622866340687
$ws.Range("H14").Value = 1.251915089846136
$ws.Range("I14").Value = 1.150391823455493
$ws.Range("J14").Value = 0.1617348924925963
$ws.Range("L14").Value = 0.3896363112056207
$ws.Range("M14").Value = 0.4325543008550738
$ws.Range("N14").Value = 1.667111867383959

$ws.Range("B15").Value = 1.841738153845995
$ws.Range("C15").Value = 0.2341158890764632
$ws.Range("D15").Value = 0.1016278879638044
$ws.Range("F15").Value = 1.993275762541117
$ws.Range("G15").Value = 1.388317638639649
$ws.Range("H15").Value = 1.251894445106814
$ws.Range("I15").Value = 1.150760197752568
$ws.Range("J15").Value = 0.1617993485165083
$ws.Range("L15").Value = 0.3891748261135035
$ws.Range("M15").Value = 0.4306324801009467
$ws.Range("N15").Value = 1.669200323401942

$ws.Range("B16").Value = 1.7739468568044
$ws.Range("C16").Value = 0.2201209332241945
$ws.Range("D16").Value = 0.101532102134442
$ws.Range("F16").Value = 1.988780106194696
$ws.Range("G16").Value = 1.381100745415154
$ws.Range("H16").Value = 1.251958207934933
$ws.Range("I16").Value = 1.153054687615636
$ws.Range("J16").Value = 0.1621750731718974
$ws.Range("L16").Value = 0.3865714787886958
$ws.Range("M16").Value = 0.4196625759026276
$ws.Range("N16").Value = 1.681352656792477

$ws.Range("B17").Value = 1.732494537092634
$ws.Range("C17").Value = 0.2115234625224218
$ws.Range("D17").Value = 0.1014762176777388
$ws.Range("F17").Value = 1.98630735502249
$ws.Range("G17").Value = 1.376907677391117
$ws.Range("H17").Value = 1.252159755708931
$ws.Range("I17").Value = 1.154625904047322
$ws.Range("J17").Value = 0.1624112465226801
$ws.Range("L17").Value = 0.3850112827880281
$ws.Range("M17").Value = 0.4129716710152778
$ws.Range("N17").Value = 1.688972046788528

$ws.Range("B18").Value = 1.708702380439206
$ws.Range("C18").Value = 0.2065740208639966
$ws.Range("D18").Value = 0.1014451351946022
$ws.Range("F18").Value = 1.984989853957572
$ws.Range("G18").Value = 1.374581917666319
$ws.Range("H18").Value = 1.252335360850566
$ws.Range("I18").Value = 1.155589714346334
$ws.Range("J18").Value = 0.1625491777082333
$ws.Range("L18").Value = 0.3841274805349713
$ws.Range("M18").Value = 0.4091375539631485
$ws.Range("N18").Value = 1.693414880685795

$ws.Range("B19").Value = 1.700655406218686
$ws.Range("C19").Value = 0.2048974698107315
$ws.Range("D19").Value = 0.1014347937197648
$ws.Range("F19").Value = 1.98456175905649
$ws.Range("G19").Value = 1.373809212625602
$ws.Range("H19").Value = 1.252405066676914
$ws.Range("I19").Value = 1.155926360181525
$ws.Range("J19").Value = 0.1625962382395265
$ws.Range("L19").Value = 0.3838305757563489
$ws.Range("M19").Value = 0.407841852453501
$ws.Range("N19").Value = 1.694929512654904

$ws.Range("B20").Value = 1.736902027727695
$ws.Range("C20").Value = 0.2124391340364582
$ws.Range("D20").Value = 0.1014820569669475
$ws.Range("F20").Value = 1.98655974030504
$ws.Range("G20").Value = 1.377345135369055
$ws.Range("H20").Value = 1.252132123638887
$ws.Range("I20").Value = 1.154452425525918
$ws.Range("J20").Value = 0.1623858891892773
$ws.Range("L20").Value = 0.3851759633982681
$ws.Range("M20").Value = 0.413682449222911
$ws.Range("N20").Value = 1.688154701832218

$ws.Range("B21").Value = 1.859282127265999
$ws.Range("C21").Value = 0.2377253854559171
$ws.Range("D21").Value = 0.1016535046397102
$ws.Range("F21").Value = 1.994524659087858
$ws.Range("G21").Value = 1.390253411774182
$ws.Range("H21").Value = 1.251928210695723
$ws.Range("I21").Value = 1.150218425194609
$ws.Range("J21").Value = 0.1617040963304852
$ws.Range("L21").Value = 0.3898583342269006
$ws.Range("M21").Value = 0.4334766293297818
$ws.Range("N21").Value = 1.666113655977973

$ws.Range("B22").Value = 1.939703590991428
$ws.Range("C22").Value = 0.2542119233870324
$ws.Range("D22").Value = 0.1017749283843976
$ws.Range("F22").Value = 2.000664046722392
$ws.Range("G22").Value = 1.39945768502136
$ws.Range("H22").Value = 1.252326662181162
$ws.Range("I22").Value = 1.147987236386122
$ws.Range("J22").Value = 0.1612770425249206
$ws.Range("L22").Value = 0.3930388285190674
$ws.Range("M22").Value = 0.4465393729259048
$ws.Range("N22").Value = 1.652246670003411

$ws.Range("B23").Value = 1.896741328115979
$ws.Range("C23").Value = 0.2454163742115156
$ws.Range("D23").Value = 0.1017092681496585
$ws.Range("F23").Value = 1.997301844684685
$ws.Range("G23").Value = 1.394474809994222
$ws.Range("H23").Value = 1.25206533503183
$ws.Range("I23").Value = 1.149128978393897
$ws.Range("J23").Value = 0.1615032768043596
$ws.Range("L23").Value = 0.3913303630297236
$ws.Range("M23").Value = 0.4395560683617603
$ws.Range("N23").Value = 1.659598474520543

$ws.Range("B24").Value = 1.734909278605528
$ws.Range("C24").Value = 0.212025179632235
$ws.Range("D24").Value = 0.1014794137652437
$ws.Range("F24").Value = 1.986445312605667
$ws.Range("G24").Value = 1.37714709620829
$ws.Range("H24").Value = 1.252144430029119
$ws.Range("I24").Value = 1.154530666724767
$ws.Range("J24").Value = 0.1623973465353297
$ws.Range("L24").Value = 0.3851014702483724
$ws.Range("M24").Value = 0.4133610672152628
$ws.Range("N24").Value = 1.688524029283219

$ws.Range("B25").Value = 1.562267159665794
$ws.Range("C25").Value = 0.1758390329826796
$ws.Range("D25").Value = 0.1012720777451683
$ws.Range("F25").Value = 1.978739289756717
$ws.Range("G25").Value = 1.361742129581543
$ws.Range("H25").Value = 1.254511249537771
$ws.Range("I25").Value = 1.162654011275521
$ws.Range("J25").Value = 0.1634417200011207
$ws.Range("L25").Value = 0.3789021912532604
$ws.Range("M25").Value = 0.3856534927106807
$ws.Range("N25").Value = 1.722032463460701
